$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) "Find the correlation of total_penalties and total_offense." item
#    -> prefix with "Using R," (kept as its own run, ahead of the
#    pre-existing leading-space run) and lower-case the initial "F".
# ---------------------------------------------------------------------
$r1 = $d.Content
$found1 = $r1.Find.Execute("Find the correlation of total_penalties and total_offense.", `
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found1) {
    # Insert "Using R," immediately before the existing leading space run
    # (i.e. right before the word "Find") so it becomes its own run.
    $insertPos1 = $d.Range($r1.Start - 1, $r1.Start - 1)
    $insertPos1.InsertBefore("Using R,")

    # Re-locate "Find the correlation..." (position shifted) and lower-case
    # just the leading "F".
    $r1b = $d.Content
    $found1b = $r1b.Find.Execute("Find the correlation of total_penalties and total_offense.", `
        $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if ($found1b) {
        $charRange1 = $d.Range($r1b.Start, $r1b.Start + 1)
        $charRange1.Text = "f"
    }
}

# ---------------------------------------------------------------------
# 2) "Test the significance ... Check your work with R." item
#    -> prefix with "Using R, " and drop the trailing
#    "Check your work with R." sentence.
# ---------------------------------------------------------------------
$oldText2 = "Test the significance of the correlation between the total_offense and the total_penalties of a player. Provide an interpretation of the results. Check your work with R."
$newText2 = "Using R, test the significance of the correlation between the total_offense and the total_penalties of a player. Provide an interpretation of the results. "
$d.Content.Find.Execute($oldText2, $true, $false, $false, $false, $false, $true, 1, $false, $newText2, 2) | Out-Null

# ---------------------------------------------------------------------
# 3) The "H0: ... Ha:" paragraph right after it becomes a proper
#    ListParagraph-styled paragraph, keeping its 360-twip left indent but
#    adding a 360-twip first-line indent as well.
# ---------------------------------------------------------------------
$targetPara = $null
$prevText = ""
foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text
    if (($t -match "^H0:\t+Ha:") -and ($prevText -match "^Using R, test the significance")) {
        $targetPara = $p
        break
    }
    $prevText = $t
}
if ($targetPara -ne $null) {
    $targetPara.Range.Select()
    $word.Selection.Style = "ListParagraph"
    $targetPara.Format.FirstLineIndent = 18
}
